$wb = $excel.ActiveWorkbook

# Sheet "展览" (worksheet 1): update "想去人数" (column F) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1833
$ws1.Range("F12").Value = 5272
$ws1.Range("F16").Value = 2345
$ws1.Range("F18").Value = 42
$ws1.Range("F19").Value = 2186

# Sheet "全部类型" (worksheet 4): update the same counts (duplicated rows)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1833
$ws4.Range("F12").Value = 5272
$ws4.Range("F18").Value = 2345
$ws4.Range("F21").Value = 42
$ws4.Range("F22").Value = 2186
